{"js": "// Apply the five text replacements described by the diff.\nconst replacements = [\n  {\n    find: \"Several incisions 1/4\\u201d long\",\n    replace: \"Several 1/4\\u201d incisions 1/4\\u201d\",\n  },\n  {\n    find: \"A telescope is used to examine the abdomen\",\n    replace: \"Telescope examines the abdomen\",\n  },\n  {\n    find: \"A PCP is critical to coordinate care between specialists.\",\n    replace: \"Critical to coordinate care between specialists.\",\n  },\n  {\n    find: \"Call our referral line at (844) 235-6998 if you need a PCP\",\n    replace: \"PCP Referral Line (844) 235-6998\",\n  },\n  {\n    find: \"Every day counts! (Aim for some activity every day)\",\n    replace: \"Every day counts! (Aim for daily activity)\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the five text replacements described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Several incisions 1/4\u201d long\"; Replace = \"Several 1/4\u201d incisions 1/4\u201d\" },\n    @{ Find = \"A telescope is used to examine the abdomen\"; Replace = \"Telescope examines the abdomen\" },\n    @{ Find = \"A PCP is critical to coordinate care between specialists.\"; Replace = \"Critical to coordinate care between specialists.\" },\n    @{ Find = \"Call our referral line at (844) 235-6998 if you need a PCP\"; Replace = \"PCP Referral Line (844) 235-6998\" },\n    @{ Find = \"Every day counts! (Aim for some activity every day)\"; Replace = \"Every day counts! (Aim for daily activity)\" }\n)\n\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text\n    foreach ($r in $replacements) {\n        if ($ptext -eq ($r.Find + \"`r\")) {\n            $p.Range.Text = $r.Replace\n            break\n        }\n    }\n}\n"}
